$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.73908466666667
$ws.Range("H2").Value = 38.217254
$ws.Range("I2").Value = 0.005953388968763418
$ws.Range("J2").Value = 0.006105597140986208
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 205.6181175590836
$ws.Range("R2").Value = 1850.563058031752
$ws.Range("S2").Value = 0.0002037873623224447
$ws.Range("T2").Value = 0.000221930419280361
$ws.Range("G3").Value = 12.73908466666667
$ws.Range("H3").Value = 38.217254
$ws.Range("I3").Value = 0.005953388968763418
$ws.Range("J3").Value = 0.006105597140986208
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 1032.672450131116
$ws.Range("R3").Value = 9294.052051180048
$ws.Range("S3").Value = 0.001023477878571696
$ws.Range("T3").Value = 0.00111459745161329
$ws.Range("G4").Value = 12.73908466666667
$ws.Range("H4").Value = 38.217254
$ws.Range("I4").Value = 0.005953388968763418
$ws.Range("J4").Value = 0.006105597140986208
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 2149.13084527078
$ws.Range("R4").Value = 19342.17760743702
$ws.Range("S4").Value = 0.002129995700002898
$ws.Range("T4").Value = 0.002319627838447889
$ws.Range("G5").Value = 12.73908466666667
$ws.Range("H5").Value = 38.217254
$ws.Range("I5").Value = 0.005953388968763418
$ws.Range("J5").Value = 0.006105597140986208
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 1050.143546354905
$ws.Range("R5").Value = 6300.861278129431
$ws.Range("S5").Value = 0.001040793417973542
$ws.Range("T5").Value = 0.0007556363881865965
$ws.Range("G6").Value = 12.73908466666667
$ws.Range("H6").Value = 38.217254
$ws.Range("I6").Value = 0.005953388968763418
$ws.Range("J6").Value = 0.006105597140986208
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 1569.307198523145
$ws.Range("R6").Value = 14123.76478670831
$ws.Range("S6").Value = 0.001555334609892838
$ws.Range("T6").Value = 0.001693805043458071
$ws.Range("I7").Value = 0.9182810852447438
$ws.Range("J7").Value = 0.9417584502053091
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 31715.58739548569
$ws.Range("R7").Value = 285440.2865593712
$ws.Range("S7").Value = 0.03143320236834585
$ws.Range("T7").Value = 0.0342316800287822
$ws.Range("I8").Value = 0.9182810852447438
$ws.Range("J8").Value = 0.9417584502053091
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.1578664491754216
$ws.Range("T8").Value = 0.1719211969600354
$ws.Range("I9").Value = 0.9182810852447438
$ws.Range("J9").Value = 0.9417584502053091
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 331492.9051810508
$ws.Range("R9").Value = 2983436.146629457
$ws.Range("S9").Value = 0.328541402758632
$ws.Range("T9").Value = 0.3577912311844601
$ws.Range("I10").Value = 0.9182810852447438
$ws.Range("J10").Value = 0.9417584502053091
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 161979.4978069184
$ws.Range("R10").Value = 971876.9868415105
$ws.Range("S10").Value = 0.1605372863064997
$ws.Range("T10").Value = 0.1165532113280571
$ws.Range("I11").Value = 0.9182810852447438
$ws.Range("J11").Value = 0.9417584502053091
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 242057.9479861446
$ws.Range("R11").Value = 2178521.531875302
$ws.Range("S11").Value = 0.2399027446358447
$ws.Range("T11").Value = 0.2612611307039743
$ws.Range("G12").Value = 1.091866333333334
$ws.Range("H12").Value = 3.275599000000001
$ws.Range("I12").Value = 0.0005102646818291153
$ws.Range("J12").Value = 0.0005233104369407934
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 17.62351895451245
$ws.Range("R12").Value = 158.6116705906121
$ws.Range("S12").Value = [double]"1.746660501133958E-05"
$ws.Range("T12").Value = [double]"1.902164554952932E-05"
$ws.Range("G13").Value = 1.091866333333334
$ws.Range("H13").Value = 3.275599000000001
$ws.Range("I13").Value = 0.0005102646818291153
$ws.Range("J13").Value = 0.0005233104369407934
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 88.51030597271679
$ws.Range("R13").Value = 796.5927537544512
$ws.Range("S13").Value = [double]"8.772223968712063E-05"
$ws.Range("T13").Value = [double]"9.553209390468086E-05"
$ws.Range("G14").Value = 1.091866333333334
$ws.Range("H14").Value = 3.275599000000001
$ws.Range("I14").Value = 0.0005102646818291153
$ws.Range("J14").Value = 0.0005233104369407934
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 184.2019012574301
$ws.Range("R14").Value = 1657.81711131687
$ws.Range("S14").Value = 0.0001825618288779668
$ws.Range("T14").Value = 0.0001988151903324102
$ws.Range("G15").Value = 1.091866333333334
$ws.Range("H15").Value = 3.275599000000001
$ws.Range("I15").Value = 0.0005102646818291153
$ws.Range("J15").Value = 0.0005233104369407934
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 90.00775278874254
$ws.Range("R15").Value = 540.0465167324552
$ws.Range("S15").Value = [double]"8.920635373542836E-05"
$ws.Range("T15").Value = [double]"6.476555844403753E-05"
$ws.Range("G16").Value = 1.091866333333334
$ws.Range("H16").Value = 3.275599000000001
$ws.Range("I16").Value = 0.0005102646818291153
$ws.Range("J16").Value = 0.0005233104369407934
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 134.5052444159179
$ws.Range("R16").Value = 1210.547199743261
$ws.Range("S16").Value = 0.0001333076545172599
$ws.Range("T16").Value = 0.0001451759487101354
$ws.Range("G17").Value = 160.0313415
$ws.Range("H17").Value = 320.062683
$ws.Range("I17").Value = 0.0747878554913321
$ws.Range("J17").Value = 0.05113328661083746
$ws.Range("M17").Value = 16.14072933333334
$ws.Range("N17").Value = 48.42218800000001
$ws.Range("O17").Value = 0.03423048004954622
$ws.Range("P17").Value = 0.03634868370049611
$ws.Range("Q17").Value = 2583.022568001734
$ws.Range("R17").Value = 15498.1354080104
$ws.Range("S17").Value = 0.002560024195344389
$ws.Range("T17").Value = 0.001858627661584144
$ws.Range("G18").Value = 160.0313415
$ws.Range("H18").Value = 320.062683
$ws.Range("I18").Value = 0.0747878554913321
$ws.Range("J18").Value = 0.05113328661083746
$ws.Range("O18").Value = 0.1719151703242873
$ws.Range("P18").Value = 0.1825533892714798
$ws.Range("Q18").Value = 12972.67125926219
$ws.Range("R18").Value = 77836.02755557316
$ws.Range("S18").Value = 0.01285716691498055
$ws.Range("T18").Value = 0.00933455477539836
$ws.Range("G19").Value = 160.0313415
$ws.Range("H19").Value = 320.062683
$ws.Range("I19").Value = 0.0747878554913321
$ws.Range("J19").Value = 0.05113328661083746
$ws.Range("M19").Value = 168.70371
$ws.Range("N19").Value = 506.11113
$ws.Range("O19").Value = 0.3577786889414888
$ws.Range("P19").Value = 0.3799182594076638
$ws.Range("Q19").Value = 26997.88102732696
$ws.Range("R19").Value = 161987.2861639618
$ws.Range("S19").Value = 0.02675750088643432
$ws.Range("T19").Value = 0.01942646924698257
$ws.Range("G20").Value = 160.0313415
$ws.Range("H20").Value = 320.062683
$ws.Range("I20").Value = 0.0747878554913321
$ws.Range("J20").Value = 0.05113328661083746
$ws.Range("M20").Value = 82.43477250000001
$ws.Range("N20").Value = 164.869545
$ws.Range("O20").Value = 0.1748236883957081
$ws.Range("P20").Value = 0.1237612588479007
$ws.Range("Q20").Value = 13192.14722942231
$ws.Range("R20").Value = 52768.58891768924
$ws.Range("S20").Value = 0.01307468874419989
$ws.Range("T20").Value = 0.00632831991998775
$ws.Range("G21").Value = 160.0313415
$ws.Range("H21").Value = 320.062683
$ws.Range("I21").Value = 0.0747878554913321
$ws.Range("J21").Value = 0.05113328661083746
$ws.Range("M21").Value = 123.1883796666667
$ws.Range("N21").Value = 369.565139
$ws.Range("O21").Value = 0.2612519722889696
$ws.Range("P21").Value = 0.2774184087724594
$ws.Range("Q21").Value = 19714.00165526799
$ws.Range("R21").Value = 118284.0099316079
$ws.Range("S21").Value = 0.01953847475037295
$ws.Range("T21").Value = 0.01418531500688463
$ws.Range("G22").Value = 1.000156333333333
$ws.Range("H22").Value = 3.000469
$ws.Range("I22").Value = 0.0004674056133315229
$ws.Range("J22").Value = 0.0004793556059265206
$ws.Range("M22").Value = 16.14072933333334
$ws.Range("N22").Value = 48.42218800000001
$ws.Range("O22").Value = 0.03423048004954622
$ws.Range("P22").Value = 0.03634868370049611
$ws.Range("Q22").Value = 16.14325266735245
$ws.Range("R22").Value = 145.289274006172
$ws.Range("S22").Value = [double]"1.599951852219061E-05"
$ws.Range("T22").Value = [double]"1.742394529988276E-05"
$ws.Range("G23").Value = 1.000156333333333
$ws.Range("H23").Value = 3.000469
$ws.Range("I23").Value = 0.0004674056133315229
$ws.Range("J23").Value = 0.0004793556059265206
$ws.Range("O23").Value = 0.1719151703242873
$ws.Range("P23").Value = 0.1825533892714798
$ws.Range("Q23").Value = 81.07598923178676
$ws.Range("R23").Value = 729.683903086081
$ws.Range("S23").Value = [double]"8.035411562641675E-05"
$ws.Range("T23").Value = [double]"8.75079905281702E-05"
$ws.Range("G24").Value = 1.000156333333333
$ws.Range("H24").Value = 3.000469
$ws.Range("I24").Value = 0.0004674056133315229
$ws.Range("J24").Value = 0.0004793556059265206
$ws.Range("M24").Value = 168.70371
$ws.Range("N24").Value = 506.11113
$ws.Range("O24").Value = 0.3577786889414888
$ws.Range("P24").Value = 0.3799182594076638
$ws.Range("Q24").Value = 168.73008401333
$ws.Range("R24").Value = 1518.57075611997
$ws.Range("S24").Value = 0.0001672277675416447
$ws.Range("T24").Value = 0.0001821159474409097
$ws.Range("G25").Value = 1.000156333333333
$ws.Range("H25").Value = 3.000469
$ws.Range("I25").Value = 0.0004674056133315229
$ws.Range("J25").Value = 0.0004793556059265206
$ws.Range("M25").Value = 82.43477250000001
$ws.Range("N25").Value = 164.869545
$ws.Range("O25").Value = 0.1748236883957081
$ws.Range("P25").Value = 0.1237612588479007
$ws.Range("Q25").Value = 82.4476598027675
$ws.Range("R25").Value = 494.685958816605
$ws.Range("S25").Value = [double]"8.171357329947497E-05"
$ws.Range("T25").Value = [double]"5.932565322526439E-05"
$ws.Range("G26").Value = 1.000156333333333
$ws.Range("H26").Value = 3.000469
$ws.Range("I26").Value = 0.0004674056133315229
$ws.Range("J26").Value = 0.0004793556059265206
$ws.Range("M26").Value = 123.1883796666667
$ws.Range("N26").Value = 369.565139
$ws.Range("O26").Value = 0.2612519722889696
$ws.Range("P26").Value = 0.2774184087724594
$ws.Range("Q26").Value = 123.2076381166879
$ws.Range("R26").Value = 1108.868743050191
$ws.Range("S26").Value = 0.0001221106383417959
$ws.Range("T26").Value = 0.0001329820694322935
